# Apply updated cryptocurrency price/volume data
# (values scraped on Sun Sep 22 19:11:51 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.932.03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.574.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.15"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.63%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.86%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.58"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.349"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.04"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.039.30"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.845.22"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.571.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.33"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.62%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.75"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.59"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.86%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.165"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.31%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.16%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "459.68"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0797"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.66"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.80"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.398"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.51"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.00%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.55%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.05"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "158.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.633"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0959"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.04"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.14%  "
